# Commit: "Added MG to medical specialty."
#
# The "medical_specialty" field's Description cell (column D) lists the two-
# letter medical specialty panel codes. Insert a new "MG = Medical Genetics"
# line right after the existing "IM = Immunology" line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the row whose "Field Name" column (B) is "medical_specialty" so the
# edit is resilient to row shuffling.
$targetRow = 0
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $fieldName = $ws.Cells.Item($r, 2).Value2
    if ($fieldName -eq "medical_specialty") {
        $targetRow = $r
        break
    }
}

if ($targetRow -eq 0) {
    # Fallback to the known row from the source workbook.
    $targetRow = 9
}

$descCell = $ws.Cells.Item($targetRow, 4)
$oldText = $descCell.Value2
$newText = $oldText -replace [regex]::Escape("IM = Immunology`n"), "IM = Immunology`nMG = Medical Genetics`n"
$descCell.Value2 = $newText
